# ============================================================
# portfolio_history_detailed.xlsx -- "new condicion small - consecultives bloks"
# Appends 8 new trade-log rows (171-178) to Sheet1 and applies a
# tiny floating-point correction to the timestamp already stored
# in A170. The used range grows from A1:S170 to A1:S178.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the existing date/time display format from column A so the
# new timestamp cells render identically to the rest of the column.
$dateFmt = $ws.Range("A170").NumberFormat

# --- Correct the stored timestamp precision on the existing last row (170) ---
$ws.Range("A170").Value = 45606.58189127315

# --- Append the new trade rows 171-178 ---
# Columns B:S hold numeric-looking data (quantities, prices, totals, ...)
# that the source log stores as literal TEXT (inlineStr), not numbers --
# e.g. "0.00" must stay "0.00", not become the number 0. Setting the
# cell to the Text number format before assigning the value keeps Excel
# from auto-converting these numeric-looking strings into real numbers.

# -- Row 171 --
$ws.Range("A171").Value = 45607.64172494213
$ws.Range("A171").NumberFormat = $dateFmt
$ws.Range("B171:S171").NumberFormat = "@"
$ws.Range("B171").Value = "BTCUSDT"
$ws.Range("C171").Value = "buy"
$ws.Range("D171").Value = "0.01229300"
$ws.Range("E171").Value = "84415.88"
$ws.Range("F171").Value = "1037.72"
$ws.Range("G171").Value = "Compra em oportunidade de curto prazo com base em indicadores e limite de compras consecutivas"
$ws.Range("H171").Value = "0.00"
$ws.Range("I171").Value = "73388.77"
$ws.Range("J171").Value = "0.12293000"
$ws.Range("K171").Value = "1037.72"
$ws.Range("L171").Value = "0.00"
$ws.Range("M171").Value = "0.00"
$ws.Range("N171").Value = "72351.05"
$ws.Range("O171").Value = "0.00"
$ws.Range("P171").Value = "0.00"
$ws.Range("Q171").Value = "0.00"
$ws.Range("R171").Value = "Neutral"
$ws.Range("S171").Value = "Loss of 0.00"

# -- Row 172 --
$ws.Range("A172").Value = 45607.64188423611
$ws.Range("A172").NumberFormat = $dateFmt
$ws.Range("B172:S172").NumberFormat = "@"
$ws.Range("B172").Value = "BTCUSDT"
$ws.Range("C172").Value = "sell"
$ws.Range("D172").Value = "0.01229300"
$ws.Range("E172").Value = "84467.24"
$ws.Range("F172").Value = "1038.36"
$ws.Range("G172").Value = "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas"
$ws.Range("H172").Value = "1038.36"
$ws.Range("I172").Value = "73389.40"
$ws.Range("J172").Value = "0.12293000"
$ws.Range("K172").Value = "0.00"
$ws.Range("L172").Value = "1038.36"
$ws.Range("M172").Value = "1038.36"
$ws.Range("N172").Value = "73389.40"
$ws.Range("O172").Value = "0.00"
$ws.Range("P172").Value = "1.41"
$ws.Range("Q172").Value = "0.00"
$ws.Range("R172").Value = "Good"
$ws.Range("S172").Value = "Profit of 1038.36"

# -- Row 173 --
$ws.Range("A173").Value = 45607.64204055555
$ws.Range("A173").NumberFormat = $dateFmt
$ws.Range("B173:S173").NumberFormat = "@"
$ws.Range("B173").Value = "BTCUSDT"
$ws.Range("C173").Value = "sell"
$ws.Range("D173").Value = "0.01229300"
$ws.Range("E173").Value = "84467.24"
$ws.Range("F173").Value = "1038.36"
$ws.Range("G173").Value = "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas"
$ws.Range("H173").Value = "1038.36"
$ws.Range("I173").Value = "73389.40"
$ws.Range("J173").Value = "0.12293000"
$ws.Range("K173").Value = "0.00"
$ws.Range("L173").Value = "1038.36"
$ws.Range("M173").Value = "2076.71"
$ws.Range("N173").Value = "73389.40"
$ws.Range("O173").Value = "0.00"
$ws.Range("P173").Value = "2.83"
$ws.Range("Q173").Value = "0.00"
$ws.Range("R173").Value = "Good"
$ws.Range("S173").Value = "Profit of 1038.36"

# -- Row 174 --
$ws.Range("A174").Value = 45607.64218829861
$ws.Range("A174").NumberFormat = $dateFmt
$ws.Range("B174:S174").NumberFormat = "@"
$ws.Range("B174").Value = "BTCUSDT"
$ws.Range("C174").Value = "sell"
$ws.Range("D174").Value = "0.01229300"
$ws.Range("E174").Value = "84467.24"
$ws.Range("F174").Value = "1038.36"
$ws.Range("G174").Value = "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas"
$ws.Range("H174").Value = "1038.36"
$ws.Range("I174").Value = "73389.40"
$ws.Range("J174").Value = "0.12293000"
$ws.Range("K174").Value = "0.00"
$ws.Range("L174").Value = "1038.36"
$ws.Range("M174").Value = "3115.07"
$ws.Range("N174").Value = "73389.40"
$ws.Range("O174").Value = "0.00"
$ws.Range("P174").Value = "4.24"
$ws.Range("Q174").Value = "0.00"
$ws.Range("R174").Value = "Good"
$ws.Range("S174").Value = "Profit of 1038.36"

# -- Row 175 --
$ws.Range("A175").Value = 45607.66492768518
$ws.Range("A175").NumberFormat = $dateFmt
$ws.Range("B175:S175").NumberFormat = "@"
$ws.Range("B175").Value = "BTCUSDT"
$ws.Range("C175").Value = "sell"
$ws.Range("D175").Value = "0.00983500"
$ws.Range("E175").Value = "85537.14"
$ws.Range("F175").Value = "841.26"
$ws.Range("G175").Value = "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas"
$ws.Range("H175").Value = "841.26"
$ws.Range("I175").Value = "75071.24"
$ws.Range("J175").Value = "0.09835000"
$ws.Range("K175").Value = "0.00"
$ws.Range("L175").Value = "841.26"
$ws.Range("M175").Value = "841.26"
$ws.Range("N175").Value = "75071.24"
$ws.Range("O175").Value = "0.00"
$ws.Range("P175").Value = "1.12"
$ws.Range("Q175").Value = "0.00"
$ws.Range("R175").Value = "Good"
$ws.Range("S175").Value = "Profit of 841.26"

# -- Row 176 --
$ws.Range("A176").Value = 45607.66507328704
$ws.Range("A176").NumberFormat = $dateFmt
$ws.Range("B176:S176").NumberFormat = "@"
$ws.Range("B176").Value = "BTCUSDT"
$ws.Range("C176").Value = "sell"
$ws.Range("D176").Value = "0.00983500"
$ws.Range("E176").Value = "85536.70"
$ws.Range("F176").Value = "841.25"
$ws.Range("G176").Value = "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas"
$ws.Range("H176").Value = "841.25"
$ws.Range("I176").Value = "75071.24"
$ws.Range("J176").Value = "0.09835000"
$ws.Range("K176").Value = "0.00"
$ws.Range("L176").Value = "841.25"
$ws.Range("M176").Value = "1682.51"
$ws.Range("N176").Value = "75071.24"
$ws.Range("O176").Value = "0.00"
$ws.Range("P176").Value = "2.24"
$ws.Range("Q176").Value = "0.00"
$ws.Range("R176").Value = "Good"
$ws.Range("S176").Value = "Profit of 841.25"

# -- Row 177 --
$ws.Range("A177").Value = 45607.72814283278
$ws.Range("A177").NumberFormat = $dateFmt
$ws.Range("B177:S177").NumberFormat = "@"
$ws.Range("B177").Value = "BTCUSDT"
$ws.Range("C177").Value = "buy"
$ws.Range("D177").Value = "0.00472500"
$ws.Range("E177").Value = "86849.99"
$ws.Range("F177").Value = "410.37"
$ws.Range("G177").Value = "Compra em oportunidade de curto prazo com base em indicadores e limite de compras consecutivas"
$ws.Range("H177").Value = "0.00"
$ws.Range("I177").Value = "79455.73"
$ws.Range("J177").Value = "0.04725000"
$ws.Range("K177").Value = "410.37"
$ws.Range("L177").Value = "0.00"
$ws.Range("M177").Value = "0.00"
$ws.Range("N177").Value = "79045.36"
$ws.Range("O177").Value = "0.00"
$ws.Range("P177").Value = "0.00"
$ws.Range("Q177").Value = "0.00"
$ws.Range("R177").Value = "Neutral"
$ws.Range("S177").Value = "Loss of 0.00"

# -- Row 178 --
$ws.Range("A178").Value = 45607.7282989224
$ws.Range("A178").NumberFormat = $dateFmt
$ws.Range("B178:S178").NumberFormat = "@"
$ws.Range("B178").Value = "BTCUSDT"
$ws.Range("C178").Value = "sell"
$ws.Range("D178").Value = "0.00472500"
$ws.Range("E178").Value = "86849.99"
$ws.Range("F178").Value = "410.37"
$ws.Range("G178").Value = "Venda para lucro a curto prazo em carteira pequena com limite de vendas consecutivas"
$ws.Range("H178").Value = "410.37"
$ws.Range("I178").Value = "79455.73"
$ws.Range("J178").Value = "0.04725000"
$ws.Range("K178").Value = "0.00"
$ws.Range("L178").Value = "410.37"
$ws.Range("M178").Value = "410.37"
$ws.Range("N178").Value = "79455.73"
$ws.Range("O178").Value = "0.00"
$ws.Range("P178").Value = "1.41"
$ws.Range("Q178").Value = "0.00"
$ws.Range("R178").Value = "Good"
$ws.Range("S178").Value = "Profit of 410.37"

# The sheet dimension (A1:S178) is recomputed automatically from the used range.

